$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell I1: new shared string "5.0.9", same style as H1 (s=3)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "5.0.9"

# Rows 2-11: empty I cells with style matching H column (s=2)
foreach ($r in 2..11) {
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
}

# Row 12: I12 date value with style matching F12 (s=5, date style)
$ws.Range("F12").Copy()
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I12").Value = 42358.649265277774

# Rows 13-15: empty I cells with style matching G column (s=2)
foreach ($r in 13..15) {
    $ws.Range("G$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
}

# Row 17: I17 date value with style matching F17 (s=5, date style)
$ws.Range("F17").Copy()
$ws.Range("I17").PasteSpecial(-4122)
$ws.Range("I17").Value = 42358.667345023146

# Update selection to I17
$ws.Range("I17").Select()

$excel.CutCopyMode = $false
